# Update "想去人数" (column F) counts on the sheets that contain them.
# Values below map worksheet row number -> new value for column F.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览" = @{
        2  = 118
        3  = 187
        4  = 433
        5  = 204
        6  = 136
        7  = 1201
        8  = 402
        9  = 201
        10 = 57
        12 = 381
        13 = 421
        14 = 796
        16 = 731
        17 = 292
        18 = 82
        19 = 1023
        20 = 475
        21 = 274
        22 = 89
        25 = 45
        26 = 476
    }
    "演出" = @{
        4  = 369
        13 = 5
    }
    "全部类型" = @{
        4  = 118
        5  = 187
        6  = 433
        7  = 204
        8  = 136
        9  = 1201
        10 = 402
        11 = 201
        13 = 57
        14 = 369
        17 = 381
        20 = 421
        21 = 796
        23 = 731
        24 = 292
        25 = 82
        26 = 1023
        27 = 475
        30 = 274
        31 = 89
        36 = 45
        38 = 476
        39 = 5
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
